$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 5).Value2 = 3
$ws.Cells.Item(2, 7).Value2 = 88.72291666666666
$ws.Cells.Item(2, 8).Value2 = 266.16875
$ws.Cells.Item(2, 9).Value2 = 0.7675060578750151
$ws.Cells.Item(2, 10).Value2 = 0.7675060578750152
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 13).Value2 = 170.93328
$ws.Cells.Item(2, 14).Value2 = 512.79984
$ws.Cells.Item(2, 15).Value2 = 0.7687311215213114
$ws.Cells.Item(2, 16).Value2 = 0.7687311215213115
$ws.Cells.Item(2, 17).Value2 = 15165.699157
$ws.Cells.Item(2, 18).Value2 = 136491.292413
$ws.Cells.Item(2, 19).Value2 = 0.5900057926446609
$ws.Cells.Item(2, 20).Value2 = 0.5900057926446611

# Row 3
$ws.Cells.Item(3, 5).Value2 = 3
$ws.Cells.Item(3, 7).Value2 = 88.72291666666666
$ws.Cells.Item(3, 8).Value2 = 266.16875
$ws.Cells.Item(3, 9).Value2 = 0.7675060578750151
$ws.Cells.Item(3, 10).Value2 = 0.7675060578750152
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 13).Value2 = 40.31217066666667
$ws.Cells.Item(3, 14).Value2 = 120.936512
$ws.Cells.Item(3, 15).Value2 = 0.1812942463137967
$ws.Cells.Item(3, 16).Value2 = 0.1812942463137967
$ws.Cells.Item(3, 17).Value2 = 3576.613358711111
$ws.Cells.Item(3, 18).Value2 = 32189.5202284
$ws.Cells.Item(3, 19).Value2 = 0.1391444323037241
$ws.Cells.Item(3, 20).Value2 = 0.1391444323037241

# Row 4
$ws.Cells.Item(4, 5).Value2 = 3
$ws.Cells.Item(4, 7).Value2 = 88.72291666666666
$ws.Cells.Item(4, 8).Value2 = 266.16875
$ws.Cells.Item(4, 9).Value2 = 0.7675060578750151
$ws.Cells.Item(4, 10).Value2 = 0.7675060578750152
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 13).Value2 = 11.112244
$ws.Cells.Item(4, 14).Value2 = 33.336732
$ws.Cells.Item(4, 15).Value2 = 0.04997463216489184
$ws.Cells.Item(4, 16).Value2 = 0.04997463216489184
$ws.Cells.Item(4, 17).Value2 = 985.9106983916665
$ws.Cells.Item(4, 18).Value2 = 8873.196285524999
$ws.Cells.Item(4, 19).Value2 = 0.03835583292663007
$ws.Cells.Item(4, 20).Value2 = 0.03835583292663008

# Row 5
$ws.Cells.Item(5, 5).Value2 = 3
$ws.Cells.Item(5, 7).Value2 = 17.91585
$ws.Cells.Item(5, 8).Value2 = 53.74755
$ws.Cells.Item(5, 9).Value2 = 0.1549827702197958
$ws.Cells.Item(5, 10).Value2 = 0.1549827702197958
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 13).Value2 = 170.93328
$ws.Cells.Item(5, 14).Value2 = 512.79984
$ws.Cells.Item(5, 15).Value2 = 0.7687311215213114
$ws.Cells.Item(5, 16).Value2 = 0.7687311215213115
$ws.Cells.Item(5, 17).Value2 = 3062.415004488
$ws.Cells.Item(5, 18).Value2 = 27561.735040392
$ws.Cells.Item(5, 19).Value2 = 0.1191400787675433
$ws.Cells.Item(5, 20).Value2 = 0.1191400787675434

# Row 6
$ws.Cells.Item(6, 5).Value2 = 3
$ws.Cells.Item(6, 7).Value2 = 17.91585
$ws.Cells.Item(6, 8).Value2 = 53.74755
$ws.Cells.Item(6, 9).Value2 = 0.1549827702197958
$ws.Cells.Item(6, 10).Value2 = 0.1549827702197958
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 13).Value2 = 40.31217066666667
$ws.Cells.Item(6, 14).Value2 = 120.936512
$ws.Cells.Item(6, 15).Value2 = 0.1812942463137967
$ws.Cells.Item(6, 16).Value2 = 0.1812942463137967
$ws.Cells.Item(6, 17).Value2 = 722.2268028384001
$ws.Cells.Item(6, 18).Value2 = 6500.0412255456
$ws.Cells.Item(6, 19).Value2 = 0.02809748451862222
$ws.Cells.Item(6, 20).Value2 = 0.02809748451862222

# Row 7
$ws.Cells.Item(7, 5).Value2 = 3
$ws.Cells.Item(7, 7).Value2 = 17.91585
$ws.Cells.Item(7, 8).Value2 = 53.74755
$ws.Cells.Item(7, 9).Value2 = 0.1549827702197958
$ws.Cells.Item(7, 10).Value2 = 0.1549827702197958
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 13).Value2 = 11.112244
$ws.Cells.Item(7, 14).Value2 = 33.336732
$ws.Cells.Item(7, 15).Value2 = 0.04997463216489184
$ws.Cells.Item(7, 16).Value2 = 0.04997463216489184
$ws.Cells.Item(7, 17).Value2 = 199.0852966674
$ws.Cells.Item(7, 18).Value2 = 1791.7676700066
$ws.Cells.Item(7, 19).Value2 = 0.007745206933630249
$ws.Cells.Item(7, 20).Value2 = 0.007745206933630251

# Row 8
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 7).Value2 = 8.960212333333333
$ws.Cells.Item(8, 8).Value2 = 26.880637
$ws.Cells.Item(8, 9).Value2 = 0.077511171905189
$ws.Cells.Item(8, 10).Value2 = 0.07751117190518901
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 13).Value2 = 170.93328
$ws.Cells.Item(8, 14).Value2 = 512.79984
$ws.Cells.Item(8, 15).Value2 = 0.7687311215213114
$ws.Cells.Item(8, 16).Value2 = 0.7687311215213115
$ws.Cells.Item(8, 17).Value2 = 1531.59848363312
$ws.Cells.Item(8, 18).Value2 = 13784.38635269808
$ws.Cells.Item(8, 19).Value2 = 0.0595852501091071
$ws.Cells.Item(8, 20).Value2 = 0.05958525010910712

# Row 9
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 7).Value2 = 8.960212333333333
$ws.Cells.Item(9, 8).Value2 = 26.880637
$ws.Cells.Item(9, 9).Value2 = 0.077511171905189
$ws.Cells.Item(9, 10).Value2 = 0.07751117190518901
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 13).Value2 = 40.31217066666667
$ws.Cells.Item(9, 14).Value2 = 120.936512
$ws.Cells.Item(9, 15).Value2 = 0.1812942463137967
$ws.Cells.Item(9, 16).Value2 = 0.1812942463137967
$ws.Cells.Item(9, 17).Value2 = 361.2056087909049
$ws.Cells.Item(9, 18).Value2 = 3250.850479118144
$ws.Cells.Item(9, 19).Value2 = 0.01405232949145037
$ws.Cells.Item(9, 20).Value2 = 0.01405232949145038

# Row 10
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 7).Value2 = 8.960212333333333
$ws.Cells.Item(10, 8).Value2 = 26.880637
$ws.Cells.Item(10, 9).Value2 = 0.077511171905189
$ws.Cells.Item(10, 10).Value2 = 0.07751117190518901
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 13).Value2 = 11.112244
$ws.Cells.Item(10, 14).Value2 = 33.336732
$ws.Cells.Item(10, 15).Value2 = 0.04997463216489184
$ws.Cells.Item(10, 16).Value2 = 0.04997463216489184
$ws.Cells.Item(10, 17).Value2 = 99.56806573980931
$ws.Cells.Item(10, 18).Value2 = 896.1125916582839
$ws.Cells.Item(10, 19).Value2 = 0.003873592304631519
$ws.Cells.Item(10, 20).Value2 = 0.00387359230463152
